# BangBaoCaoTienDoCongViec.xlsx update
# - add progress note for "Tuan 3" result column (C4)
# - add new task description for "Tuan 4" (B5), which reuses the existing
#   merged-cell layout / alignment already applied to the B5:B16 column
# - move the active selection to C5:C6 (where the user was last working)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Công việc của tuần 4 (B5 đang trống, gộp ô B5:B6)
$ws.Range("B5").Value = "Vẽ Sequence Diagram, Class Diagram, Data flow Diagram(bỏ qua nếu dùng EF(code-first))"

# Kết quả của tuần 3 (A4/B4 đã có sẵn nội dung công việc)
$ws.Range("C4").Value = "Edit lại usecase Đăng nhập và thống kê"

# Excel nhớ lại vùng chọn cuối cùng của người dùng
$null = $ws.Range("C5:C6").Select()

$wb.Save()
